$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at position 3 (shifts old row3->4, old row4->5,
#    and auto-relocates the A3:A4 merge to A4:A5).
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Clear()
$ws.Range("B3").Clear()

# 2. Header row (row 1).
$ws.Range("B1").Value = "AccountType"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Account, Balance"
$ws.Range("C1").ClearContents()

# 3. Row 2 - move the numeric "Balance" value from C2 to D2, and populate
#    the new test columns.
$ws.Range("D2").Value = $ws.Range("C2").Value2
$ws.Range("C2").Value = "empty"
$ws.Range("E2").Value = "Commas ignored"

# 4. Row 3 (brand-new row) - only C3/E3 are populated.
$ws.Range("C3").Value = "header"
$ws.Range("E3").Value = "Empty rows ignored"

# 5. Row 4 (previously row 3, "bar").
$ws.Range("D4").Value = $ws.Range("C4").Value2
$ws.Range("C4").Value = "columns"
$ws.Range("E4").Value = "No header skipped"

# 6. Row 5 (previously row 4, merged blank A cell).
$ws.Range("D5").Value = $ws.Range("C5").Value2
$ws.Range("C5").Value = "ignored"

# 7. Row 6 (brand-new divider row with a thin top border, no values).
$ws.Range("A6:D6").Borders.Item(8).LineStyle = 1
$ws.Range("A6:D6").Borders.Item(8).Weight = 2

# 8. Column widths: column C matches columns A/B, column D is sized for
#    its header text.
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = 13.83

# 9. Selection moves to B1.
$ws.Range("B1").Select()
